# Simulated Wild Card round and logged it
# Updates cumulative season stats on the "Rushing" and "Receiving" sheets
# for the Panthers roster, and adds a new receiver (C.Saunders) who
# caught his first two passes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# S.Darnold (row 4): 1DATT, 2DATT, 3DATT, RZATT
$rushing.Range("C4").Value = 5
$rushing.Range("D4").Value = 1
$rushing.Range("E4").Value = 4
$rushing.Range("F4").Value = 1

# C.Hubbard (row 5)
$rushing.Range("C5").Value = 49
$rushing.Range("D5").Value = 27

# A.Abdullah (row 7)
$rushing.Range("C7").Value = 19
$rushing.Range("D7").Value = 21
$rushing.Range("F7").Value = 7

# Dj.Moore (row 10)
$rushing.Range("E10").Value = 2

# T.Tremble (row 12)
$rushing.Range("C12").Value = 1
$rushing.Range("F12").Value = 2

# ---------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# C.Hubbard (row 2): Short Target, Short Comp, RZ Target, RZ Comp
$receiving.Range("C2").Value = 23
$receiving.Range("D2").Value = 15
$receiving.Range("G2").Value = 1
$receiving.Range("H2").Value = 1

# A.Abdullah (row 4)
$receiving.Range("C4").Value = 43
$receiving.Range("D4").Value = 32
$receiving.Range("G4").Value = 7
$receiving.Range("H4").Value = 5

# R.Anderson (row 6)
$receiving.Range("C6").Value = 86
$receiving.Range("D6").Value = 49
$receiving.Range("E6").Value = 22
$receiving.Range("G6").Value = 7
$receiving.Range("H6").Value = 3

# Dj.Moore (row 7)
$receiving.Range("C7").Value = 121
$receiving.Range("D7").Value = 77
$receiving.Range("E7").Value = 42
$receiving.Range("F7").Value = 16
$receiving.Range("G7").Value = 14

# B.Zylstra (row 9)
$receiving.Range("C9").Value = 18
$receiving.Range("E9").Value = 7

# K.Kirkwood (row 11)
$receiving.Range("C11").Value = 5

# A new player, C.Saunders, recorded his first receiving stats in the
# Wild Card game. Insert a row for him right after A.Erickson (row 13)
# so the roster stays grouped the way the source sheet had it, which
# pushes T.Tremble / I.Thomas / C.Thompson down by one row.
$receiving.Range("A13:H13").Copy()
$receiving.Range("A17:H17").PasteSpecial(-4122)
$receiving.Application.CutCopyMode = $false

$receiving.Cells.Item(17, 1).Value = 15
$receiving.Cells.Item(17, 2).Value = "C.Thompson"
$receiving.Cells.Item(17, 3).Value = 1
$receiving.Cells.Item(17, 4).Value = 0
$receiving.Cells.Item(17, 5).Value = 0
$receiving.Cells.Item(17, 6).Value = 0
$receiving.Cells.Item(17, 7).Value = 0
$receiving.Cells.Item(17, 8).Value = 0

$receiving.Cells.Item(16, 1).Value = 14
$receiving.Cells.Item(16, 2).Value = "I.Thomas"
$receiving.Cells.Item(16, 3).Value = 26
$receiving.Cells.Item(16, 4).Value = 15
$receiving.Cells.Item(16, 5).Value = 4
$receiving.Cells.Item(16, 6).Value = 3
$receiving.Cells.Item(16, 7).Value = 4
$receiving.Cells.Item(16, 8).Value = 2

$receiving.Cells.Item(15, 1).Value = 13
$receiving.Cells.Item(15, 2).Value = "T.Tremble"
$receiving.Cells.Item(15, 3).Value = 28
$receiving.Cells.Item(15, 4).Value = 19
$receiving.Cells.Item(15, 5).Value = 6
$receiving.Cells.Item(15, 6).Value = 1
$receiving.Cells.Item(15, 7).Value = 4
$receiving.Cells.Item(15, 8).Value = 4

$receiving.Cells.Item(14, 1).Value = 12
$receiving.Cells.Item(14, 2).Value = "C.Saunders"
$receiving.Cells.Item(14, 3).Value = 2
$receiving.Cells.Item(14, 4).Value = 2
$receiving.Cells.Item(14, 5).Value = 0
$receiving.Cells.Item(14, 6).Value = 0
$receiving.Cells.Item(14, 7).Value = 0
$receiving.Cells.Item(14, 8).Value = 0
